# Update crypto price (D) and volume-change (E) columns per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.706.44"
$ws.Range("E2").Value = "  -0.04%  "

$ws.Range("D3").Value = "3.235.12"
$ws.Range("E3").Value = "  +0.81%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "610.05"
$ws.Range("E5").Value = "  +1.70%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.46"
$ws.Range("E6").Value = "  +2.28%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.11%  "

$ws.Range("D8").Value = "3.233.56"
$ws.Range("E8").Value = "  +0.79%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.552"
$ws.Range("E9").Value = "  +1.08%  "

$ws.Range("E10").Value = "  +0.48%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.72"
$ws.Range("E11").Value = "  -4.89%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.507"
$ws.Range("E12").Value = "  -2.46%  "

$ws.Range("E13").Value = "  +1.81%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "39.06"
$ws.Range("E14").Value = "  -0.64%  "

$ws.Range("D15").Value = "3.770.07"
$ws.Range("E15").Value = "  +0.93%  "

$ws.Range("D16").Value = "66.801.20"
$ws.Range("E16").Value = "  +0.24%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.40"
$ws.Range("E17").Value = "  -0.85%  "

$ws.Range("D18").Value = "3.237.41"
$ws.Range("E18").Value = "  +0.86%  "

$ws.Range("E19").Value = "  +1.18%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "511.75"
$ws.Range("E20").Value = "  -0.80%  "

$ws.Range("E21").Value = "  -0.89%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.738"
$ws.Range("E22").Value = "  -0.86%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.03"
$ws.Range("E23").Value = "  -0.92%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.68"
$ws.Range("E24").Value = "  -3.17%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.26"
$ws.Range("E25").Value = "  -0.47%  "

$ws.Range("E26").Value = "  +0.17%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.01"
$ws.Range("E27").Value = "  -0.51%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.15"
$ws.Range("E28").Value = "  -1.54%  "

$ws.Range("E29").Value = "  +2.09%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.123"
$ws.Range("E30").Value = "  +36.36%  "

$ws.Range("E31").Value = "  +0.89%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.01"
$ws.Range("E32").Value = "  -2.67%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "28.25"
$ws.Range("E33").Value = "  -0.48%  "

$ws.Range("E34").Value = "  +0.10%  "

$ws.Range("E35").Value = "  -3.66%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.52"
$ws.Range("E36").Value = "  -0.80%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "508.96"
$ws.Range("E37").Value = "  +2.29%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "55.53"
$ws.Range("E38").Value = "  +1.24%  "

$ws.Range("D39").Value = "0.0₃0772"
$ws.Range("E39").Value = "  +13.00%  "

$ws.Range("E40").Value = "  +7.13%  "

$ws.Range("E41").Value = "  -0.09%  "

$ws.Range("E42").Value = "  +6.28%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.77"
$ws.Range("E43").Value = "  -1.64%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.300"
$ws.Range("E44").Value = "  -0.63%  "

$ws.Range("E45").Value = "  +0.33%  "

$ws.Range("D46").Value = "2.919.41"
$ws.Range("E46").Value = "  -0.60%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.23"
$ws.Range("E47").Value = "  -1.71%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.44"
$ws.Range("E48").Value = "  +4.06%  "

$ws.Range("E49").Value = "  -0.07%  "

$ws.Range("E50").Value = "  -0.63%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "123.32"
$ws.Range("E51").Value = "  +0.50%  "

